$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# 1. Fill in the "network number" (H) column for the two existing
#    AOC-based MBUF rows (46 and 47), which now reference the new network.
$ws.Range("H46").Value = "NGF_Networks_NGFround2_P2_01"
$ws.Range("H47").Value = "NGF_Networks_NGFround2_P2_01"

# 2. Insert two new rows before row 54 (shifts the existing P4/P5/P6
#    calibration rows that used to start at row 54 down to start at row 56).
$ws.Rows("54:55").Insert()

# 3. Populate the two newly inserted rows (54 and 55) with the two new
#    AOC-based MBUF runs added for toll-cap post-processing.
$ws.Range("A54").Value = "NextGenFwys"
$ws.Range("B54").Value = 2035
$ws.Range("C54").Value = "2035_TM160_NGF_r2_NoProject_06_add2.9cT"
$ws.Range("D54").Value = "NGF_Round2"
$ws.Range("E54").Value = "NoProject"
$ws.Range("F54").Value = "aoc based mbuf at 5c, no cap, network with more transit"
$ws.Range("G54").Value = "current"
$ws.Range("H54").Value = "NGF_Networks_NGFround2_P2_01"
$ws.Range("I54").Value = "https://app.asana.com/0/1203644633064654/1208050271846456/f"

$ws.Range("A55").Value = "NextGenFwys"
$ws.Range("B55").Value = 2035
$ws.Range("C55").Value = "2035_TM160_NGF_r2_NoProject_06_add4.4cT"
$ws.Range("D55").Value = "NGF_Round2"
$ws.Range("E55").Value = "NoProject"
$ws.Range("F55").Value = "aoc based mbuf at 7.5c, no cap, network with more transit"
$ws.Range("G55").Value = "current"
$ws.Range("H55").Value = "NGF_Networks_NGFround2_P2_01"
$ws.Range("I55").Value = "https://app.asana.com/0/1203644633064654/1208050271846456/f"

# 4. The data range grew by 2 rows (A1:L65 -> A1:L67); keep the
#    _FilterDatabase defined name in sync with the new extent.
$wb.Names.Item("all_runs!_FilterDatabase").RefersTo = "=all_runs!`$A`$1:`$L`$67"
